# Auto-generated Excel COM-interop script applying market price updates
# to the Leve-profit tables across all 8 job sheets (scheduled runner update).
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 266.66666
$ws.Range("I2").Value = 280
$ws.Range("K2").Value = 280
$ws.Range("M2").Value = -167
$ws.Range("H33").Value = 751.5
$ws.Range("I33").Value = 751.5
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 751.5
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -522.5
$ws.Range("N33").ClearContents()
$ws.Range("H62").Value = 5381.3076
$ws.Range("I62").Value = 5444.087
$ws.Range("K62").Value = 5444.087
$ws.Range("M62").Value = -4820.087
$ws.Range("H65").Value = 5381.3076
$ws.Range("I65").Value = 5444.087
$ws.Range("K65").Value = 27220.435
$ws.Range("M65").Value = -24100.435
$ws.Range("H98").Value = 7240
$ws.Range("I98").Value = 7277.5
$ws.Range("J98").Value = 7202.5
$ws.Range("K98").Value = 7277.5
$ws.Range("L98").Value = 7202.5
$ws.Range("M98").Value = -5779.5
$ws.Range("N98").Value = -10198.5
$ws.Range("H112").Value = 30575.514
$ws.Range("I112").Value = 1276.75
$ws.Range("J112").Value = 38657.93
$ws.Range("K112").Value = 3830.25
$ws.Range("L112").Value = 115973.79
$ws.Range("M112").Value = -2722.25
$ws.Range("N112").Value = -118189.79
$ws.Range("H122").Value = 7240
$ws.Range("I122").Value = 7277.5
$ws.Range("J122").Value = 7202.5
$ws.Range("K122").Value = 21832.5
$ws.Range("L122").Value = 21607.5
$ws.Range("M122").Value = -19382.5
$ws.Range("N122").Value = -26507.5
$ws.Range("H133").Value = 59148
$ws.Range("J133").Value = 59148
$ws.Range("L133").Value = 59148
$ws.Range("N133").Value = -69268
$ws.Range("H136").Value = 94993.336
$ws.Range("J136").Value = 94993.336
$ws.Range("L136").Value = 94993.336
$ws.Range("N136").Value = -105193.336
$ws.Range("H138").Value = 5321797
$ws.Range("I138").Value = 1171.7778
$ws.Range("J138").Value = 7465929
$ws.Range("K138").Value = 3515.3334
$ws.Range("L138").Value = 22397787
$ws.Range("M138").Value = 1624.6666
$ws.Range("N138").Value = -22408067
$ws.Range("H139").Value = 83850
$ws.Range("J139").Value = 83850
$ws.Range("L139").Value = 83850
$ws.Range("N139").Value = -94130
$ws.Range("H140").Value = 87374.75
$ws.Range("J140").Value = 87374.75
$ws.Range("L140").Value = 87374.75
$ws.Range("N140").Value = -97734.75

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6021.2607
$ws.Range("I32").Value = 3630.5156
$ws.Range("K32").Value = 3630.5156
$ws.Range("M32").Value = -3343.5156
$ws.Range("H61").Value = 3212.675
$ws.Range("I61").Value = 2451.8386
$ws.Range("K61").Value = 2451.8386
$ws.Range("M61").Value = -2239.8386
$ws.Range("H63").Value = 3201.4285
$ws.Range("I63").Value = 3250
$ws.Range("J63").Value = 3182
$ws.Range("K63").Value = 3250
$ws.Range("L63").Value = 3182
$ws.Range("M63").Value = -2564
$ws.Range("N63").Value = -4554
$ws.Range("H66").Value = 3201.4285
$ws.Range("I66").Value = 3250
$ws.Range("J66").Value = 3182
$ws.Range("K66").Value = 16250
$ws.Range("L66").Value = 15910
$ws.Range("M66").Value = -12818
$ws.Range("N66").Value = -22774
$ws.Range("H110").Value = 9162.3125
$ws.Range("I110").Value = 10363.363
$ws.Range("J110").Value = 6520
$ws.Range("K110").Value = 10363.363
$ws.Range("L110").Value = 6520
$ws.Range("M110").Value = -8318.362999999999
$ws.Range("N110").Value = -10610
$ws.Range("H132").Value = 3178.4595
$ws.Range("I132").Value = 3045.8
$ws.Range("K132").Value = 9137.400000000001
$ws.Range("M132").Value = -6607.400000000001
$ws.Range("H136").Value = 3212.675
$ws.Range("I136").Value = 2451.8386
$ws.Range("K136").Value = 7355.5158
$ws.Range("M136").Value = -4805.5158

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1819.4359
$ws.Range("I134").Value = 1817.091
$ws.Range("K134").Value = 5451.272999999999
$ws.Range("M134").Value = -2916.272999999999

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4597.4443
$ws.Range("I16").Value = 4561.857
$ws.Range("J16").Value = 4722
$ws.Range("K16").Value = 4561.857
$ws.Range("L16").Value = 4722
$ws.Range("M16").Value = -4274.857
$ws.Range("N16").Value = -5296
$ws.Range("H31").Value = 70310.92999999999
$ws.Range("I31").Value = 94081.37
$ws.Range("K31").Value = 94081.37
$ws.Range("M31").Value = -93786.37
$ws.Range("H34").Value = 70310.92999999999
$ws.Range("I34").Value = 94081.37
$ws.Range("K34").Value = 94081.37
$ws.Range("M34").Value = -93879.37
$ws.Range("H58").Value = 2512.0208
$ws.Range("I58").Value = 2171.8718
$ws.Range("K58").Value = 2171.8718
$ws.Range("M58").Value = -1968.8718
$ws.Range("H113").Value = 4597.4443
$ws.Range("I113").Value = 4561.857
$ws.Range("J113").Value = 4722
$ws.Range("K113").Value = 4561.857
$ws.Range("L113").Value = 4722
$ws.Range("M113").Value = -2391.857
$ws.Range("N113").Value = -9062
$ws.Range("H132").Value = 4764.5557
$ws.Range("I132").Value = 4983.143
$ws.Range("J132").Value = 3999.5
$ws.Range("K132").Value = 14949.429
$ws.Range("L132").Value = 11998.5
$ws.Range("M132").Value = -12419.429
$ws.Range("N132").Value = -17058.5
$ws.Range("H134").Value = 11912.077
$ws.Range("I134").Value = 6788.56
$ws.Range("K134").Value = 20365.68
$ws.Range("M134").Value = -17830.68
$ws.Range("H136").Value = 2512.0208
$ws.Range("I136").Value = 2171.8718
$ws.Range("K136").Value = 6515.6154
$ws.Range("M136").Value = -3965.6154

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 2500
$ws.Range("I64").Value = 1750
$ws.Range("J64").Value = 4000
$ws.Range("K64").Value = 5250
$ws.Range("L64").Value = 12000
$ws.Range("M64").Value = -4980
$ws.Range("N64").Value = -12540
$ws.Range("H67").Value = 2500
$ws.Range("I67").Value = 1750
$ws.Range("J67").Value = 4000
$ws.Range("K67").Value = 5250
$ws.Range("L67").Value = 12000
$ws.Range("M67").Value = -4314
$ws.Range("N67").Value = -13872
$ws.Range("H107").Value = 538.9375
$ws.Range("I107").Value = 452.45456
$ws.Range("J107").Value = 729.2
$ws.Range("K107").Value = 1357.36368
$ws.Range("L107").Value = 2187.6
$ws.Range("M107").Value = 562.6363200000001
$ws.Range("N107").Value = -6027.6
$ws.Range("H110").Value = 9499.75
$ws.Range("I110").Value = 9499.75
$ws.Range("K110").Value = 28499.25
$ws.Range("M110").Value = -24409.25
$ws.Range("H119").Value = 1532.6666
$ws.Range("I119").Value = 1532.6666
$ws.Range("K119").Value = 4597.9998
$ws.Range("M119").Value = 240.0002000000004
$ws.Range("H139").Value = 2148.5652
$ws.Range("I139").Value = 1969.3684
$ws.Range("J139").Value = 2999.75
$ws.Range("K139").Value = 5908.1052
$ws.Range("L139").Value = 8999.25
$ws.Range("M139").Value = -768.1052
$ws.Range("N139").Value = -19279.25
$ws.Range("H140").Value = 4470
$ws.Range("I140").Value = 3778.75
$ws.Range("K140").Value = 11336.25
$ws.Range("M140").Value = -6156.25

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5846
$ws.Range("I70").Value = 6186.5713
$ws.Range("J70").Value = 5250
$ws.Range("K70").Value = 6186.5713
$ws.Range("L70").Value = 5250
$ws.Range("M70").Value = -5916.5713
$ws.Range("N70").Value = -5790
$ws.Range("H73").Value = 5846
$ws.Range("I73").Value = 6186.5713
$ws.Range("J73").Value = 5250
$ws.Range("K73").Value = 6186.5713
$ws.Range("L73").Value = 5250
$ws.Range("M73").Value = -5250.5713
$ws.Range("N73").Value = -7122
$ws.Range("H122").Value = 4487.7617
$ws.Range("I122").Value = 3888.75
$ws.Range("J122").Value = 6404.6
$ws.Range("K122").Value = 11666.25
$ws.Range("L122").Value = 19213.8
$ws.Range("M122").Value = -9216.25
$ws.Range("N122").Value = -24113.8
$ws.Range("H132").Value = 3990.0715
$ws.Range("I132").Value = 3208.3
$ws.Range("K132").Value = 9624.900000000001
$ws.Range("M132").Value = -7094.900000000001

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6472.875
$ws.Range("I7").Value = 6677.5557
$ws.Range("J7").Value = 6209.7144
$ws.Range("K7").Value = 6677.5557
$ws.Range("L7").Value = 6209.7144
$ws.Range("M7").Value = -6565.5557
$ws.Range("N7").Value = -6433.7144
$ws.Range("H61").Value = 4189.8667
$ws.Range("I61").Value = 3765.5833
$ws.Range("K61").Value = 3765.5833
$ws.Range("M61").Value = -3563.5833
$ws.Range("H113").Value = 4189.8667
$ws.Range("I113").Value = 3765.5833
$ws.Range("K113").Value = 3765.5833
$ws.Range("M113").Value = -1595.5833
$ws.Range("H126").Value = 6472.875
$ws.Range("I126").Value = 6677.5557
$ws.Range("J126").Value = 6209.7144
$ws.Range("K126").Value = 20032.6671
$ws.Range("L126").Value = 18629.1432
$ws.Range("M126").Value = -17562.6671
$ws.Range("N126").Value = -23569.1432
$ws.Range("H132").Value = 3203.2666
$ws.Range("I132").Value = 2515.125
$ws.Range("J132").Value = 5955.8335
$ws.Range("K132").Value = 7545.375
$ws.Range("L132").Value = 17867.5005
$ws.Range("M132").Value = -5015.375
$ws.Range("N132").Value = -22927.5005
$ws.Range("H136").Value = 5737.5
$ws.Range("I136").Value = 4566.6665
$ws.Range("J136").Value = 9250
$ws.Range("K136").Value = 13699.9995
$ws.Range("L136").Value = 27750
$ws.Range("M136").Value = -11149.9995
$ws.Range("N136").Value = -32850

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1691.3704
$ws.Range("I126").Value = 1503.1364
$ws.Range("J126").Value = 2519.6
$ws.Range("K126").Value = 4509.4092
$ws.Range("L126").Value = 7558.799999999999
$ws.Range("M126").Value = -2039.4092
$ws.Range("N126").Value = -12498.8
$ws.Range("H132").Value = 2258.3948
$ws.Range("I132").Value = 1949.6
$ws.Range("J132").Value = 3416.375
$ws.Range("K132").Value = 5848.799999999999
$ws.Range("L132").Value = 10249.125
$ws.Range("M132").Value = -3318.799999999999
$ws.Range("N132").Value = -15309.125
$ws.Range("H136").Value = 1611.9459
$ws.Range("I136").Value = 1407.862
$ws.Range("J136").Value = 2351.75
$ws.Range("K136").Value = 4223.586
$ws.Range("L136").Value = 7055.25
$ws.Range("M136").Value = -1673.586
$ws.Range("N136").Value = -12155.25
